# Add five more three-character words (with readings) below the existing
# "こうか" row, matching the author's upload of additional rows to the
# words_3_chars workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (A3) was blank; it now holds the word that used to trail the list.
$ws.Range("A3").Value = "たいか"

# New rows 4-7, each a single word in column A (column B left blank).
$ws.Range("A4").Value = "けいか"
$ws.Range("A5").Value = "そうか"
$ws.Range("A6").Value = "さいか"
$ws.Range("A7").Value = "せいか"

# Mirror the saved selection/window state from the authored workbook.
$ws.Range("F15").Select()
